$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.962.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.859.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '697.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.44%  '
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.856.18'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.29%  '
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.37%  '
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.20'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.15%  '
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.63%  '
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.72%  '
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.31'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.510.65'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.858.54'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.37%  '
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '71.018.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.82%  '
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '499.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.14%  '
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.79'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.41%  '
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.721'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.59%  '
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000150'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.43%  '
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.03'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.27%  '
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.77%  '
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.42%  '
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.79%  '
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.15'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.51%  '
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.53'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.27'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.90%  '
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.56'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.182'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.34%  '
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.813.88'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.46%  '
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.55%  '
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.35%  '
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +8.67%  '
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.17%  '
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.60%  '
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000312'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.79%  '
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '163.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.73%  '
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '49.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '416.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.70%  '
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.68%  '
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.20%  '
$ws.Range("E51").Style = "Normal"
